$wb = $excel.ActiveWorkbook

# --- Sheet "a1a" ---
$ws = $wb.Worksheets.Item("a1a")
$ws.Range("B7").Value = 0.356189
$ws.Range("C7").Value = 0.286787
$ws.Range("D7").Value = 0.011347
$ws.Range("E7").Value = 0.168465
$ws.Range("G7").Value = 0.197524
$ws.Range("H7").Value = 0.758722

$ws.Range("B8").Value = 0.596524
$ws.Range("C8").Value = 0.508024
$ws.Range("D8").Value = 0.00155
$ws.Range("E8").Value = 0.336445
$ws.Range("G8").Value = 0.19507
$ws.Range("H8").Value = 0.758948

$ws.Range("B9").Value = 1.124542
$ws.Range("C9").Value = 0.963901
$ws.Range("D9").Value = 0.002528
$ws.Range("E9").Value = 0.645697
$ws.Range("G9").Value = 0.187301
$ws.Range("H9").Value = 0.768284

$ws.Range("H9").Select()

# --- Sheet "ijcnn1" ---
$ws = $wb.Worksheets.Item("ijcnn1")
$ws.Range("B7").Value = 9.257949
$ws.Range("C7").Value = 4.242995
$ws.Range("D7").Value = 9.194281
$ws.Range("E7").Value = 0.025517
$ws.Range("G7").Value = 0.609846
$ws.Range("H7").Value = 0.904996

$ws.Range("B8").Value = 5.450975
$ws.Range("C8").Value = 2.388517
$ws.Range("D8").Value = 5.254436
$ws.Range("E8").Value = 0.097093
$ws.Range("G8").Value = 0.642672
$ws.Range("H8").Value = 0.904996

$ws.Range("B9").Value = 3.231666
$ws.Range("C9").Value = 1.467572
$ws.Range("D9").Value = 2.709832
$ws.Range("E9").Value = 1.117482
$ws.Range("G9").Value = 0.627083
$ws.Range("H9").Value = 0.917329

$ws.Range("H9").Select()

# --- Sheet "generated" ---
$ws = $wb.Worksheets.Item("generated")
$ws.Range("B7").Value = 59.094832
$ws.Range("C7").Value = 38.845353
$ws.Range("D7").Value = 36.273127
$ws.Range("E7").Value = 0.326354
$ws.Range("G7").Value = 0.614886
$ws.Range("H7").Value = 0.6011

$ws.Range("B8").Value = 42.116462
$ws.Range("C8").Value = 20.788727
$ws.Range("D8").Value = 38.433629
$ws.Range("G8").Value = 0.643389
$ws.Range("H8").Value = 0.6011

$ws.Range("B9").Value = 44.992385
$ws.Range("C9").Value = 23.649346
$ws.Range("D9").Value = 38.27289
$ws.Range("E9").Value = 13.576266
$ws.Range("G9").Value = 0.647126
$ws.Range("H9").Value = 0.6016

$ws.Range("G9").Select()
